# Apply header renames and title-case fixes for Spanish connector words (de/del/el/la/las/los/y)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellUpdates = @(
  @("A1", "mx_state"),
  @("B1", "mx_municipality"),
  @("C1", "n_matriculas"),
  @("D1", "pct_matriculas"),
  @("B6", "Pabellón De Arteaga"),
  @("B7", "Rincón De Romos"),
  @("B8", "San Francisco De Los Romo"),
  @("B30", "Amatenango De La Frontera"),
  @("B32", "Bejucal De Ocampo"),
  @("B40", "Comitán De Domínguez"),
  @("B57", "Mazapa De Madero"),
  @("B67", "Salto De Agua"),
  @("B68", "San Cristóbal De Las Casas"),
  @("B98", "Hidalgo Del Parral"),
  @("B105", "San Francisco Del Oro"),
  @("B119", "San Juan De Sabinas"),
  @("A127", "Ciudad De México"),
  @("B131", "Cuajimalpa De Morelos"),
  @("B153", "Nombre De Dios"),
  @("B155", "Pánuco De Coronado"),
  @("B157", "San Juan Del Río"),
  @("B158", "San Pedro Del Gallo"),
  @("A164", "Estado De México"),
  @("B164", "Acambay De Ruíz Castañeda"),
  @("B167", "Almoloya De Alquisiras"),
  @("B168", "Almoloya De Juárez"),
  @("B171", "Atizapán De Zaragoza"),
  @("B176", "Coacalco De Berriozábal"),
  @("B181", "Ecatepec De Morelos"),
  @("B183", "Ixtapan De La Sal"),
  @("B184", "Ixtapan Del Oro"),
  @("B195", "Naucalpan De Juárez"),
  @("B202", "San Felipe Del Progreso"),
  @("B211", "Tenango Del Valle"),
  @("B216", "Tlalnepantla De Baz"),
  @("B222", "Valle De Bravo"),
  @("B231", "San Miguel De Allende"),
  @("B232", "Apaseo El Alto"),
  @("B233", "Apaseo El Grande"),
  @("B241", "Dolores Hidalgo Cuna De La Independencia Nacional"),
  @("B245", "Jaral Del Progreso"),
  @("B253", "Purísima Del Rincón"),
  @("B257", "San Diego De La Unión"),
  @("B259", "San Francisco Del Rincón"),
  @("B261", "San Luis De La Paz"),
  @("B262", "Santa Cruz De Juventino Rosas"),
  @("B264", "Silao De La Victoria"),
  @("B269", "Valle De Santiago"),
  @("B274", "Acapulco De Juárez"),
  @("B276", "Ajuchitlán Del Progreso"),
  @("B279", "Atenango Del Río"),
  @("B280", "Atoyac De Álvarez"),
  @("B281", "Ayutla De Los Libres"),
  @("B283", "Chilapa De Álvarez"),
  @("B287", "Coyuca De Benítez"),
  @("B288", "Coyuca De Catalán"),
  @("B291", "Cutzamala De Pinzón"),
  @("B293", "Huitzuco De Los Figueroa"),
  @("B294", "Iguala De La Independencia"),
  @("B296", "Zihuatanejo De Azueta"),
  @("B307", "Taxco De Alarcón"),
  @("B309", "Técpan De Galeana"),
  @("B311", "Tixtla De Guerrero"),
  @("B314", "Tlalixtaquilla De Maldonado"),
  @("B315", "Tlapa De Comonfort"),
  @("B326", "Atotonilco El Grande"),
  @("B329", "Cuautepec De Hinojosa"),
  @("B331", "Huasca De Ocampo"),
  @("B334", "Huejutla De Reyes"),
  @("B337", "Jacala De Ledezma"),
  @("B342", "Mixquiahuala De Juárez"),
  @("B343", "Molango De Escamilla"),
  @("B345", "Omitlán De Juárez"),
  @("B346", "Pachuca De Soto"),
  @("B349", "Progreso De Obregón"),
  @("B355", "Tenango De Doria"),
  @("B357", "Tepehuacán De Guerrero"),
  @("B358", "Tepeji Del Río De Ocampo"),
  @("B360", "Tezontepec De Aldama"),
  @("B366", "Tulancingo De Bravo"),
  @("B368", "Zacualtipán De Ángeles"),
  @("B373", "Acatlán De Juárez"),
  @("B374", "Ahualulco De Mercado"),
  @("B377", "Atemajac De Brizuela"),
  @("B379", "Atotonilco El Alto"),
  @("B387", "Encarnación De Díaz"),
  @("B390", "Ixtlahuacán Del Río"),
  @("B392", "Jilotlán De Los Dolores"),
  @("B396", "Lagos De Moreno"),
  @("B402", "San Cristóbal De La Barranca"),
  @("B403", "San Diego De Alejandría"),
  @("B405", "San Juan De Los Lagos"),
  @("B409", "Tamazula De Gordiano"),
  @("B411", "Techaluta De Montenegro"),
  @("B413", "Tepatitlán De Morelos"),
  @("B416", "Tizapán El Alto"),
  @("B417", "Tlajomulco De Zúñiga"),
  @("B424", "Unión De San Antonio"),
  @("B429", "Zapotitlán De Vadillo"),
  @("B430", "Zapotlán El Grande"),
  @("B485", "Tiquicheo De Nicolás Romero"),
  @("B502", "Coatlán Del Río"),
  @("B512", "Puente De Ixtla"),
  @("B517", "Tlaltizapán De Zapata"),
  @("B526", "Ixtlán Del Río"),
  @("B531", "Santa María Del Oro"),
  @("B548", "Mier Y Noriega"),
  @("B551", "San Nicolás De Los Garza"),
  @("B554", "Acatlán De Pérez Figueroa"),
  @("B557", "Ayoquezco De Aldama"),
  @("B559", "Chalcatongo De Hidalgo"),
  @("B561", "Coicoyán De Las Flores"),
  @("B562", "Constancia Del Rosario"),
  @("B564", "Eloxochitlán De Flores Magón"),
  @("B565", "Guadalupe De Ramírez"),
  @("B566", "Heroica Ciudad De Ejutla De Crespo"),
  @("B567", "Heroica Ciudad De Huajuapan De León"),
  @("B568", "Heroica Ciudad De Tlaxiaco"),
  @("B569", "Ixtlán De Juárez"),
  @("B570", "Heroica Ciudad De Juchitán De Zaragoza"),
  @("B574", "Mariscala De Juárez"),
  @("B575", "Mártires De Tacubaya"),
  @("B578", "Miahuatlán De Porfirio Díaz"),
  @("B579", "Nejapa De Madero"),
  @("B580", "Oaxaca De Juárez"),
  @("B581", "Ocotlán De Morelos"),
  @("B582", "Pinotepa De Don Luis"),
  @("B583", "Putla Villa De Guerrero"),
  @("B584", "Reforma De Pineda"),
  @("B592", "San Antonino El Alto"),
  @("B597", "San Felipe Jalapa De Díaz"),
  @("B599", "San Francisco Del Mar"),
  @("B609", "San Juan De Los Cués"),
  @("B642", "San Pedro El Alto"),
  @("B659", "Santa Ana Del Valle"),
  @("B689", "Santiago Del Río"),
  @("B707", "Santo Domingo De Morelos"),
  @("B719", "Tamazulápam Del Espíritu Santo"),
  @("B720", "Tataltepec De Valdés"),
  @("B721", "Teotitlán De Flores Magón"),
  @("B722", "Tepelmeme Villa De Morelos"),
  @("B723", "Heroica Villa Tezoatlán De Segura Y Luna, Cuna De La Independencia De Oaxaca"),
  @("B724", "Tlacolula De Matamoros"),
  @("B725", "Totontepec Villa De Morelos"),
  @("B727", "Villa De Chilapa De Díaz"),
  @("B728", "Villa De Tututepec"),
  @("B731", "Villa Sola De Vega"),
  @("B732", "Zapotitlán Del Río"),
  @("B734", "Zimatlán De Álvarez"),
  @("B745", "Ayotoxco De Guerrero"),
  @("B765", "Izúcar De Matamoros"),
  @("B772", "Los Reyes De Juárez"),
  @("B776", "Palmar De Bravo"),
  @("B786", "San Salvador El Seco"),
  @("B787", "San Salvador El Verde"),
  @("B795", "Tepanco De López"),
  @("B796", "Tepatlaxco De Hidalgo"),
  @("B800", "Tetela De Ocampo"),
  @("B804", "Tlacotepec De Benito Juárez"),
  @("B823", "Amealco De Bonfil"),
  @("B825", "Cadereyta De Montes"),
  @("B831", "Jalpan De Serra"),
  @("B832", "Landa De Matamoros"),
  @("B835", "Pinal De Amoles"),
  @("B838", "San Juan Del Río"),
  @("B854", "Ciudad Del Maíz"),
  @("B864", "San Ciro De Acosta"),
  @("B868", "Santa María Del Río"),
  @("B875", "Villa De Guadalupe"),
  @("B876", "Villa De La Paz"),
  @("B877", "Villa De Ramos"),
  @("B911", "Jalpa De Méndez"),
  @("B939", "Soto La Marina"),
  @("B953", "Muñoz De Domingo Arenas"),
  @("B954", "Nanacamilpa De Mariano Arista"),
  @("B960", "Tetla De La Solidaridad"),
  @("B973", "Alto Lucero De Gutiérrez Barrios"),
  @("B976", "Amatlán De Los Reyes"),
  @("B988", "Castillo De Teayo"),
  @("B1002", "Cosamaloapan De Carpio"),
  @("B1003", "Cosautlán De Carvajal"),
  @("B1017", "Hueyapan De Ocampo"),
  @("B1018", "Huiloapan De Cuauhtémoc"),
  @("B1019", "Ignacio De La Llave"),
  @("B1021", "Ixhuatlán De Madero"),
  @("B1022", "Ixhuatlán Del Café"),
  @("B1029", "Juchique De Ferrer"),
  @("B1031", "Lerdo De Tejada"),
  @("B1034", "Martínez De La Torre"),
  @("B1047", "Paso De Ovejas"),
  @("B1048", "Paso Del Macho"),
  @("B1050", "Poza Rica De Hidalgo"),
  @("B1054", "Sayula De Alemán"),
  @("B1055", "Soledad De Doblado"),
  @("B1085", "Vega De Alatorre"),
  @("B1098", "Concepción Del Oro"),
  @("B1106", "Jiménez Del Teul"),
  @("B1110", "Mezquital Del Oro"),
  @("B1112", "Nochistlán De Mejía"),
  @("B1113", "Noria De Ángeles"),
  @("B1121", "Teúl De González Ortega"),
  @("B1124", "Villa De Cos")
)

foreach ($upd in $cellUpdates) {
    $ws.Range($upd[0]).Value = $upd[1]
}

# Fix two percentage values affected by floating point recalculation
$ws.Range("D133").Value = 0.009142659996549941
$ws.Range("D135").Value = 0.009142659996549941

# Remove trailing footnote/metadata rows (1130:1134)
$ws.Rows("1130:1134").Delete()
